# Updates cryptos list price (D) and volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.913.15'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '2.325.38'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.73'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.61'
$ws.Range("E6").Value = '  +1.27%  '
$ws.Range("E7").Value = '  +0.94%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.69'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  +6.37%  '
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = '2.687.07'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '2.326.03'
$ws.Range("E16").Value = '  +2.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.793'
$ws.Range("E17").Value = '  +3.16%  '
$ws.Range("D18").Value = '42.846.15'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.38'
$ws.Range("E19").Value = '  -3.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("E20").Value = '  +3.78%  '
$ws.Range("D21").Value = '0.0₃0894'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.07'
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.27'
$ws.Range("E23").Value = '  +6.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.56'
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.47'
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.45'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.51'
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.79'
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0703'
$ws.Range("E36").Value = '  +3.16%  '
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.80'
$ws.Range("E38").Value = '  +4.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.100'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.75'
$ws.Range("E40").Value = '  +3.54%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.74'
$ws.Range("E42").Value = '  +14.55%  '
$ws.Range("D43").Value = '1.937.58'
$ws.Range("E43").Value = '  -2.50%  '
$ws.Range("E44").Value = '  +1.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.33'
$ws.Range("E45").Value = '  +2.41%  '
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("D48").Value = '2.554.88'
$ws.Range("E48").Value = '  +1.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.51'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.80'
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.21'
$ws.Range("E51").Value = '  +2.56%  '
